# Fruta / hortaliza, semanal
# Update Hortaliza - Alcachofa (Terminal La Palmera de La Serena) weekly data.
# The underlying records (rows 2-7) are re-aligned: each row's Fecha (D),
# Variedad (H), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Unidad de comercializacion (N),
# Precio $/Kg (P) and Kg o Unidades (Q) are updated to reflect the
# latest weekly report, while Mercado ID/Mercado/Region/Codreg/Categoria
# ID/Categoria/Calidad/Origen/Clasificacion stay unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44426
$ws.Range("H2").Value = "Española"
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 11500
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11750
$ws.Range("N2").Value = "$/caja 30 unidades"
$ws.Range("P2").Value = 392
$ws.Range("Q2").Value = 30

# Row 3
$ws.Range("D3").Value = 44426
$ws.Range("H3").Value = "Madrigal"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 12500
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12750
$ws.Range("N3").Value = "$/caja 40 unidades"
$ws.Range("P3").Value = 319
$ws.Range("Q3").Value = 40

# Row 4
$ws.Range("D4").Value = 44420
$ws.Range("H4").Value = "Madrigal"
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("N4").Value = "$/caja 40 unidades"
$ws.Range("P4").Value = 362
$ws.Range("Q4").Value = 40

# Row 5
$ws.Range("D5").Value = 44420
$ws.Range("H5").Value = "Madrigal"
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13500
$ws.Range("N5").Value = "$/caja 40 unidades"
$ws.Range("P5").Value = 338
$ws.Range("Q5").Value = 40

# Row 6
$ws.Range("D6").Value = 44427
$ws.Range("H6").Value = "Madrigal"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("N6").Value = "$/caja 40 unidades"
$ws.Range("P6").Value = 312
$ws.Range("Q6").Value = 40

# Row 7
$ws.Range("D7").Value = 44438
$ws.Range("H7").Value = "Española"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11500
$ws.Range("N7").Value = "$/caja 30 unidades"
$ws.Range("P7").Value = 383
$ws.Range("Q7").Value = 30
